{"js": "// Replace the two-digit multiplication problems with their new values.\n// Each old problem string is unique within the document, so we can\n// safely search-and-replace each one independently.\nconst replacements = [\n  [\"15\u00d719=\", \"56\u00d754=\"],\n  [\"50\u00d779=\", \"74\u00d784=\"],\n  [\"57\u00d775=\", \"84\u00d725=\"],\n  [\"91\u00d744=\", \"28\u00d792=\"],\n  [\"56\u00d724=\", \"95\u00d712=\"],\n  [\"86\u00d773=\", \"48\u00d756=\"],\n  [\"39\u00d758=\", \"74\u00d743=\"],\n  [\"71\u00d777=\", \"68\u00d741=\"],\n  [\"55\u00d734=\", \"70\u00d717=\"],\n  [\"69\u00d765=\", \"95\u00d799=\"],\n  [\"16\u00d788=\", \"17\u00d730=\"],\n  [\"22\u00d721=\", \"25\u00d733=\"],\n  [\"32\u00d798=\", \"91\u00d772=\"],\n  [\"48\u00d768=\", \"41\u00d773=\"],\n  [\"66\u00d715=\", \"34\u00d711=\"],\n  [\"98\u00d716=\", \"53\u00d765=\"],\n  [\"94\u00d761=\", \"27\u00d763=\"],\n  [\"66\u00d772=\", \"96\u00d723=\"],\n  [\"68\u00d727=\", \"19\u00d788=\"],\n  [\"45\u00d715=\", \"37\u00d771=\"],\n  [\"69\u00d799=\", \"61\u00d770=\"],\n  [\"51\u00d715=\", \"93\u00d728=\"],\n  [\"35\u00d795=\", \"11\u00d751=\"],\n  [\"54\u00d713=\", \"55\u00d714=\"],\n  [\"58\u00d749=\", \"89\u00d787=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem with its new value.\n# Every \"old\" problem string is unique in the document, so a simple\n# Find/Replace (ReplaceAll) for each pair is sufficient and safe.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"15\u00d719=\", \"56\u00d754=\"),\n    @(\"50\u00d779=\", \"74\u00d784=\"),\n    @(\"57\u00d775=\", \"84\u00d725=\"),\n    @(\"91\u00d744=\", \"28\u00d792=\"),\n    @(\"56\u00d724=\", \"95\u00d712=\"),\n    @(\"86\u00d773=\", \"48\u00d756=\"),\n    @(\"39\u00d758=\", \"74\u00d743=\"),\n    @(\"71\u00d777=\", \"68\u00d741=\"),\n    @(\"55\u00d734=\", \"70\u00d717=\"),\n    @(\"69\u00d765=\", \"95\u00d799=\"),\n    @(\"16\u00d788=\", \"17\u00d730=\"),\n    @(\"22\u00d721=\", \"25\u00d733=\"),\n    @(\"32\u00d798=\", \"91\u00d772=\"),\n    @(\"48\u00d768=\", \"41\u00d773=\"),\n    @(\"66\u00d715=\", \"34\u00d711=\"),\n    @(\"98\u00d716=\", \"53\u00d765=\"),\n    @(\"94\u00d761=\", \"27\u00d763=\"),\n    @(\"66\u00d772=\", \"96\u00d723=\"),\n    @(\"68\u00d727=\", \"19\u00d788=\"),\n    @(\"45\u00d715=\", \"37\u00d771=\"),\n    @(\"69\u00d799=\", \"61\u00d770=\"),\n    @(\"51\u00d715=\", \"93\u00d728=\"),\n    @(\"35\u00d795=\", \"11\u00d751=\"),\n    @(\"54\u00d713=\", \"55\u00d714=\"),\n    @(\"58\u00d749=\", \"89\u00d787=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll) | Out-Null\n}\n"}
